$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export rolls its date window forward: drop the two oldest
# days (2025-08-21, 2025-08-22) and append the two newest days
# (2025-11-19, 2025-11-20) with fresh (not-yet-populated) counts of 0.
# Deleting the top two data rows shifts every remaining (date, count)
# pair up by two rows, which reproduces the "after" values exactly.
$ws.Range("A2:A3").EntireRow.Delete()

# Append the two new trailing rows. Force them to be stored as text
# (matching how every other date in column A is stored) instead of
# letting Excel auto-convert the date-shaped string into a date serial.
$ws.Cells.Item(90, 1).NumberFormat = "@"
$ws.Cells.Item(90, 1).Value = "2025-11-19"
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(90, 3).Value = 0

$ws.Cells.Item(91, 1).NumberFormat = "@"
$ws.Cells.Item(91, 1).Value = "2025-11-20"
$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(91, 3).Value = 0
